# Applies a rotation of species-observation data among rows 8, 9, 11, 12, 13, 14.
# Row 10 is left untouched. Columns affected: A, B, D, E, F, G, H, Q, R.
# The mapping (destination row <= source row, using the ORIGINAL/before values) is:
#   8 <= 14, 9 <= 12, 11 <= 8, 12 <= 13, 13 <= 11, 14 <= 9

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the "before" values for the rows involved, per relevant column.
$rows = @(8, 9, 11, 12, 13, 14)
$cols = @("A", "B", "D", "E", "F", "G", "H", "Q", "R")

$before = @{}
foreach ($r in $rows) {
    $rowData = @{}
    foreach ($c in $cols) {
        $rowData[$c] = $ws.Range("$c$r").Value2
    }
    $before[$r] = $rowData
}

# Mapping of destination row -> source row (values copied from source's "before" state)
$mapping = @{
    8  = 14
    9  = 12
    11 = 8
    12 = 13
    13 = 11
    14 = 9
}

foreach ($dst in $mapping.Keys) {
    $src = $mapping[$dst]
    $srcData = $before[$src]
    foreach ($c in $cols) {
        $ws.Range("$c$dst").Value2 = $srcData[$c]
    }
}
